$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 80
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 34
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 2450
$ws.Range("I6").Value = 2663.6365
$ws.Range("K6").Value = 7990.9095
$ws.Range("M6").Value = -7878.9095
$ws.Range("H17").Value = 1251297
$ws.Range("J17").Value = 1429910.9
$ws.Range("L17").Value = 4289732.699999999
$ws.Range("N17").Value = -4290068.699999999
$ws.Range("H19").Value = 1659.1875
$ws.Range("I19").Value = 395.5
$ws.Range("K19").Value = 395.5
$ws.Range("M19").Value = -220.5
$ws.Range("H28").Value = 1040
$ws.Range("I28").Value = 934.5
$ws.Range("K28").Value = 934.5
$ws.Range("M28").Value = -449.5
$ws.Range("H32").Value = 3163.2104
$ws.Range("J32").Value = 3422.182
$ws.Range("L32").Value = 3422.182
$ws.Range("N32").Value = -4074.182
$ws.Range("H33").Value = 5848265.5
$ws.Range("J33").Value = 1175
$ws.Range("L33").Value = 1175
$ws.Range("N33").Value = -1633
$ws.Range("H39").Value = 996.06665
$ws.Range("I39").Value = 703.5
$ws.Range("K39").Value = 2110.5
$ws.Range("M39").Value = -1814.5
$ws.Range("H40").Value = 18755936
$ws.Range("J40").Value = 25004584
$ws.Range("L40").Value = 25004584
$ws.Range("N40").Value = -25004934
$ws.Range("H105").Value = 10671
$ws.Range("J105").Value = 10671
$ws.Range("L105").Value = 10671
$ws.Range("N105").Value = -17659
$ws.Range("H111").Value = 4667.6665
$ws.Range("I111").Value = 2251.5
$ws.Range("J111").Value = 9500
$ws.Range("K111").Value = 6754.5
$ws.Range("L111").Value = 28500
$ws.Range("M111").Value = -3687.5
$ws.Range("N111").Value = -34634
$ws.Range("H116").Value = 7258.1177
$ws.Range("I116").Value = 9514.571
$ws.Range("J116").Value = 5678.6
$ws.Range("K116").Value = 9514.571
$ws.Range("L116").Value = 5678.6
$ws.Range("M116").Value = -6072.571
$ws.Range("N116").Value = -12562.6
$ws.Range("H118").Value = 76533040
$ws.Range("I118").Value = 178571730
$ws.Range("J118").Value = 4026.5
$ws.Range("K118").Value = 535715190
$ws.Range("L118").Value = 12079.5
$ws.Range("M118").Value = -535713533
$ws.Range("N118").Value = -15393.5
$ws.Range("H129").Value = 3428.7144
$ws.Range("I129").Value = 2497.5
$ws.Range("J129").Value = 3526.7368
$ws.Range("K129").Value = 7492.5
$ws.Range("L129").Value = 10580.2104
$ws.Range("M129").Value = -2492.5
$ws.Range("N129").Value = -20580.2104
$ws.Range("H132").Value = 2543.2693
$ws.Range("I132").Value = 1960.9048
$ws.Range("J132").Value = 4989.2
$ws.Range("K132").Value = 5882.7144
$ws.Range("L132").Value = 14967.6
$ws.Range("M132").Value = -3352.7144
$ws.Range("N132").Value = -20027.6
$ws.Range("H137").Value = 7787.2104
$ws.Range("I137").Value = 9180.071
$ws.Range("J137").Value = 3887.2
$ws.Range("K137").Value = 27540.213
$ws.Range("L137").Value = 11661.6
$ws.Range("M137").Value = -24990.213
$ws.Range("N137").Value = -16761.6
$ws.Range("H138").Value = 4530.07
$ws.Range("J138").Value = 4811.643
$ws.Range("L138").Value = 14434.929
$ws.Range("N138").Value = -24714.929
$ws.Range("H141").Value = 2717.3684
$ws.Range("I141").Value = 2789.6
$ws.Range("J141").Value = 2446.5
$ws.Range("K141").Value = 8368.799999999999
$ws.Range("L141").Value = 7339.5
$ws.Range("M141").Value = -3188.799999999999
$ws.Range("N141").Value = -17699.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 17529.555
$ws.Range("I28").Value = 8966.571
$ws.Range("J28").Value = 47500
$ws.Range("K28").Value = 8966.571
$ws.Range("L28").Value = 47500
$ws.Range("M28").Value = -8774.571
$ws.Range("N28").Value = -47884
$ws.Range("H32").Value = 15959.596
$ws.Range("I32").Value = 15473.956
$ws.Range("K32").Value = 15473.956
$ws.Range("M32").Value = -15186.956
$ws.Range("H45").Value = 2593.7273
$ws.Range("J45").Value = 3186.9375
$ws.Range("L45").Value = 3186.9375
$ws.Range("N45").Value = -3940.9375
$ws.Range("H61").Value = 3403.2
$ws.Range("I61").Value = 2060.3235
$ws.Range("K61").Value = 2060.3235
$ws.Range("M61").Value = -1848.3235
$ws.Range("H74").Value = 4676.058
$ws.Range("I74").Value = 4661.409
$ws.Range("K74").Value = 4661.409
$ws.Range("M74").Value = -3787.409
$ws.Range("H77").Value = 4676.058
$ws.Range("I77").Value = 4661.409
$ws.Range("K77").Value = 23307.045
$ws.Range("M77").Value = -18939.045
$ws.Range("H97").Value = 989.5789
$ws.Range("I97").Value = 700.5
$ws.Range("K97").Value = 700.5
$ws.Range("M97").Value = -204.5
$ws.Range("H99").Value = 17529.555
$ws.Range("I99").Value = 8966.571
$ws.Range("J99").Value = 47500
$ws.Range("K99").Value = 8966.571
$ws.Range("L99").Value = 47500
$ws.Range("M99").Value = -5971.571
$ws.Range("N99").Value = -53490
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H119").Value = 67000
$ws.Range("J119").Value = 67000
$ws.Range("L119").Value = 67000
$ws.Range("N119").Value = -76676
$ws.Range("H122").Value = 8002.62
$ws.Range("I122").Value = 5820.8647
$ws.Range("J122").Value = 14212.23
$ws.Range("K122").Value = 17462.5941
$ws.Range("L122").Value = 42636.69
$ws.Range("M122").Value = -15012.5941
$ws.Range("N122").Value = -47536.69
$ws.Range("H132").Value = 2292.5117
$ws.Range("I132").Value = 1771.2821
$ws.Range("K132").Value = 5313.846299999999
$ws.Range("M132").Value = -2783.846299999999
$ws.Range("H136").Value = 3403.2
$ws.Range("I136").Value = 2060.3235
$ws.Range("K136").Value = 6180.970499999999
$ws.Range("M136").Value = -3630.970499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 13373
$ws.Range("J21").Value = 13373
$ws.Range("L21").Value = 13373
$ws.Range("N21").Value = -13845
$ws.Range("H28").Value = 40520.5
$ws.Range("J28").Value = 40520.5
$ws.Range("L28").Value = 40520.5
$ws.Range("N28").Value = -41108.5
$ws.Range("H86").Value = 4199.4
$ws.Range("I86").Value = 3749.25
$ws.Range("K86").Value = 3749.25
$ws.Range("M86").Value = -2626.25
$ws.Range("H89").Value = 4199.4
$ws.Range("I89").Value = 3749.25
$ws.Range("K89").Value = 18746.25
$ws.Range("M89").Value = -13130.25
$ws.Range("H97").Value = 6278.5713
$ws.Range("I97").Value = 6278.5713
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 6278.5713
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -5287.5713
$ws.Range("N97").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 6338.8667
$ws.Range("I99").Value = 8869.799999999999
$ws.Range("K99").Value = 8869.799999999999
$ws.Range("M99").Value = -7371.799999999999
$ws.Range("H134").Value = 2561.1287
$ws.Range("I134").Value = 2706.5
$ws.Range("K134").Value = 8119.5
$ws.Range("M134").Value = -5584.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 36.333332
$ws.Range("I7").Value = 36
$ws.Range("K7").Value = 36
$ws.Range("M7").Value = 77
$ws.Range("H39").Value = 21999.75
$ws.Range("I39").Value = 12666.333
$ws.Range("K39").Value = 12666.333
$ws.Range("M39").Value = -12275.333
$ws.Range("H43").Value = 27642.625
$ws.Range("J43").Value = 27642.625
$ws.Range("L43").Value = 27642.625
$ws.Range("N43").Value = -28010.625
$ws.Range("H49").Value = 21999.75
$ws.Range("I49").Value = 12666.333
$ws.Range("K49").Value = 12666.333
$ws.Range("M49").Value = -12484.333
$ws.Range("H58").Value = 1942.25
$ws.Range("I58").Value = 1798.2858
$ws.Range("K58").Value = 1798.2858
$ws.Range("M58").Value = -1595.2858
$ws.Range("H62").Value = 17071.273
$ws.Range("I62").Value = 22248.25
$ws.Range("K62").Value = 22248.25
$ws.Range("M62").Value = -21624.25
$ws.Range("H65").Value = 17071.273
$ws.Range("I65").Value = 22248.25
$ws.Range("K65").Value = 111241.25
$ws.Range("M65").Value = -108121.25
$ws.Range("H68").Value = 184000
$ws.Range("J68").Value = 184000
$ws.Range("L68").Value = 184000
$ws.Range("N68").Value = -185498
$ws.Range("H71").Value = 184000
$ws.Range("J71").Value = 184000
$ws.Range("L71").Value = 552000
$ws.Range("N71").Value = -559488
$ws.Range("H95").Value = 36166.668
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H96").Value = 17847.4
$ws.Range("I96").Value = 6500
$ws.Range("J96").Value = 20684.25
$ws.Range("K96").Value = 6500
$ws.Range("L96").Value = 20684.25
$ws.Range("M96").Value = -3754
$ws.Range("N96").Value = -26176.25
$ws.Range("H101").Value = 27642.625
$ws.Range("J101").Value = 27642.625
$ws.Range("L101").Value = 27642.625
$ws.Range("N101").Value = -34132.625
$ws.Range("H107").Value = 4865.5713
$ws.Range("I107").Value = 654.1
$ws.Range("J107").Value = 6550.16
$ws.Range("K107").Value = 654.1
$ws.Range("L107").Value = 6550.16
$ws.Range("M107").Value = 1265.9
$ws.Range("N107").Value = -10390.16
$ws.Range("H109").Value = 45334.4
$ws.Range("J109").Value = 45279.25
$ws.Range("L109").Value = 45279.25
$ws.Range("N109").Value = -47359.25
$ws.Range("H136").Value = 1942.25
$ws.Range("I136").Value = 1798.2858
$ws.Range("K136").Value = 5394.857400000001
$ws.Range("M136").Value = -2844.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9091355
$ws.Range("J2").Value = 12500326
$ws.Range("L2").Value = 75001956
$ws.Range("N2").Value = -75002182
$ws.Range("H6").Value = 411
$ws.Range("I6").Value = 48.333332
$ws.Range("K6").Value = 144.999996
$ws.Range("M6").Value = -31.99999600000001
$ws.Range("H23").Value = 1180
$ws.Range("I23").Value = 109.666664
$ws.Range("J23").Value = 1581.375
$ws.Range("K23").Value = 328.999992
$ws.Range("L23").Value = 4744.125
$ws.Range("M23").Value = -93.99999200000002
$ws.Range("N23").Value = -5214.125
$ws.Range("H34").Value = 24999.25
$ws.Range("I34").Value = 2998.5
$ws.Range("K34").Value = 8995.5
$ws.Range("M34").Value = -8911.5
$ws.Range("H38").Value = 83.3
$ws.Range("J38").Value = 73.5
$ws.Range("L38").Value = 220.5
$ws.Range("N38").Value = -914.5
$ws.Range("H39").Value = 4299
$ws.Range("J39").Value = 5332
$ws.Range("L39").Value = 15996
$ws.Range("N39").Value = -16584
$ws.Range("H55").Value = 8090.8887
$ws.Range("J55").Value = 9831.143
$ws.Range("L55").Value = 29493.429
$ws.Range("N55").Value = -29847.429
$ws.Range("H68").Value = 10417523
$ws.Range("J68").Value = 1100
$ws.Range("L68").Value = 3300
$ws.Range("N68").Value = -4922
$ws.Range("H71").Value = 10417523
$ws.Range("J71").Value = 1100
$ws.Range("L71").Value = 9900
$ws.Range("N71").Value = -18012
$ws.Range("H113").Value = 2236.9092
$ws.Range("J113").Value = 2690.6667
$ws.Range("L113").Value = 8072.000100000001
$ws.Range("N113").Value = -12412.0001
$ws.Range("H131").Value = 13159516
$ws.Range("J131").Value = 1806
$ws.Range("L131").Value = 5418
$ws.Range("N131").Value = -15498
$ws.Range("H132").Value = 4780.6055
$ws.Range("I132").Value = 5916.154
$ws.Range("J132").Value = 2320.25
$ws.Range("K132").Value = 53245.38600000001
$ws.Range("L132").Value = 20882.25
$ws.Range("M132").Value = -50715.38600000001
$ws.Range("N132").Value = -25942.25
$ws.Range("H134").Value = 2869.875
$ws.Range("I134").Value = 2779.8572
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 8339.571599999999
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -3269.571599999999
$ws.Range("N134").Value = -20640
$ws.Range("H139").Value = 1989.3636
$ws.Range("I139").Value = 1989.3636
$ws.Range("K139").Value = 5968.0908
$ws.Range("M139").Value = -828.0907999999999
$ws.Range("H140").Value = 6411846
$ws.Range("I140").Value = 17858234
$ws.Range("K140").Value = 53574702
$ws.Range("M140").Value = -53569522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 17999
$ws.Range("J24").Value = 17999
$ws.Range("L24").Value = 17999
$ws.Range("N24").Value = -18345
$ws.Range("H80").Value = 2663.75
$ws.Range("I80").Value = 2973
$ws.Range("J80").Value = 1736
$ws.Range("K80").Value = 2973
$ws.Range("L80").Value = 1736
$ws.Range("M80").Value = -1975
$ws.Range("N80").Value = -3732
$ws.Range("H83").Value = 2663.75
$ws.Range("I83").Value = 2973
$ws.Range("J83").Value = 1736
$ws.Range("K83").Value = 14865
$ws.Range("L83").Value = 8680
$ws.Range("M83").Value = -9873
$ws.Range("N83").Value = -18664
$ws.Range("H97").Value = 1848.2858
$ws.Range("I97").Value = 1626.5454
$ws.Range("J97").Value = 2661.3333
$ws.Range("K97").Value = 1626.5454
$ws.Range("L97").Value = 2661.3333
$ws.Range("M97").Value = -1130.5454
$ws.Range("N97").Value = -3653.3333
$ws.Range("H99").Value = 31891.8
$ws.Range("I99").Value = 31891.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 31891.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -29645.8
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 9230.044
$ws.Range("I107").Value = 6453
$ws.Range("J107").Value = 17098.334
$ws.Range("K107").Value = 6453
$ws.Range("L107").Value = 17098.334
$ws.Range("M107").Value = -4533
$ws.Range("N107").Value = -20938.334
$ws.Range("H122").Value = 4937.25
$ws.Range("I122").Value = 4099.6
$ws.Range("J122").Value = 6333.3335
$ws.Range("K122").Value = 12298.8
$ws.Range("L122").Value = 19000.0005
$ws.Range("M122").Value = -9848.800000000001
$ws.Range("N122").Value = -23900.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1609.9333
$ws.Range("I22").Value = 426
$ws.Range("J22").Value = 1905.9166
$ws.Range("K22").Value = 426
$ws.Range("L22").Value = 1905.9166
$ws.Range("M22").Value = -131
$ws.Range("N22").Value = -2495.9166
$ws.Range("H27").Value = 1609.9333
$ws.Range("I27").Value = 426
$ws.Range("J27").Value = 1905.9166
$ws.Range("K27").Value = 426
$ws.Range("L27").Value = 1905.9166
$ws.Range("M27").Value = -319
$ws.Range("N27").Value = -2119.9166
$ws.Range("H40").Value = 603689.3
$ws.Range("I40").Value = 682521.2
$ws.Range("K40").Value = 682521.2
$ws.Range("M40").Value = -682385.2
$ws.Range("H46").Value = 2543
$ws.Range("I46").Value = 1399.0769
$ws.Range("K46").Value = 1399.0769
$ws.Range("M46").Value = -1211.0769
$ws.Range("H55").Value = 595.2857
$ws.Range("I55").Value = 361.25
$ws.Range("K55").Value = 361.25
$ws.Range("M55").Value = -188.25
$ws.Range("H61").Value = 2114.8333
$ws.Range("I61").Value = 2122.1765
$ws.Range("K61").Value = 2122.1765
$ws.Range("M61").Value = -1920.1765
$ws.Range("H68").Value = 2999.3
$ws.Range("I68").Value = 2999.3
$ws.Range("K68").Value = 2999.3
$ws.Range("M68").Value = -2250.3
$ws.Range("H71").Value = 2999.3
$ws.Range("I71").Value = 2999.3
$ws.Range("K71").Value = 14996.5
$ws.Range("M71").Value = -11252.5
$ws.Range("H102").Value = 59780
$ws.Range("J102").Value = 59780
$ws.Range("L102").Value = 59780
$ws.Range("N102").Value = -66270
$ws.Range("H113").Value = 2114.8333
$ws.Range("I113").Value = 2122.1765
$ws.Range("K113").Value = 2122.1765
$ws.Range("M113").Value = 47.82349999999997
$ws.Range("H122").Value = 8566.777
$ws.Range("I122").Value = 8716.166999999999
$ws.Range("K122").Value = 26148.501
$ws.Range("M122").Value = -23698.501
$ws.Range("H132").Value = 7678.3335
$ws.Range("I132").Value = 2978.4546
$ws.Range("J132").Value = 9746.280000000001
$ws.Range("K132").Value = 8935.363799999999
$ws.Range("L132").Value = 29238.84
$ws.Range("M132").Value = -6405.363799999999
$ws.Range("N132").Value = -34298.84
$ws.Range("H136").Value = 2551.5386
$ws.Range("I136").Value = 2301.087
$ws.Range("J136").Value = 4471.6665
$ws.Range("K136").Value = 6903.261
$ws.Range("L136").Value = 13414.9995
$ws.Range("M136").Value = -4353.261
$ws.Range("N136").Value = -18514.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 18624.625
$ws.Range("I51").Value = 13199.8
$ws.Range("J51").Value = 27666
$ws.Range("K51").Value = 13199.8
$ws.Range("L51").Value = 27666
$ws.Range("M51").Value = -12689.8
$ws.Range("N51").Value = -28686
$ws.Range("H52").Value = 31498
$ws.Range("I52").Value = 30000
$ws.Range("J52").Value = 31997.334
$ws.Range("K52").Value = 30000
$ws.Range("L52").Value = 31997.334
$ws.Range("M52").Value = -29774
$ws.Range("N52").Value = -32449.334
$ws.Range("H61").Value = 11500
$ws.Range("I61").Value = 11500
$ws.Range("K61").Value = 11500
$ws.Range("M61").Value = -11208
$ws.Range("H62").Value = 13168515
$ws.Range("J62").Value = 17868064
$ws.Range("L62").Value = 17868064
$ws.Range("N62").Value = -17869312
$ws.Range("H65").Value = 13168515
$ws.Range("J65").Value = 17868064
$ws.Range("L65").Value = 89340320
$ws.Range("N65").Value = -89346560
$ws.Range("H74").Value = 14315.4
$ws.Range("I74").Value = 19000
$ws.Range("K74").Value = 19000
$ws.Range("M74").Value = -18064
$ws.Range("H77").Value = 14315.4
$ws.Range("I77").Value = 19000
$ws.Range("K77").Value = 57000
$ws.Range("M77").Value = -52320
$ws.Range("H107").Value = 461.8889
$ws.Range("I107").Value = 481.31818
$ws.Range("J107").Value = 376.4
$ws.Range("K107").Value = 1443.95454
$ws.Range("L107").Value = 1129.2
$ws.Range("M107").Value = 476.04546
$ws.Range("N107").Value = -4969.2
$ws.Range("H136").Value = 3491.0315
$ws.Range("I136").Value = 3015.875
$ws.Range("K136").Value = 9047.625
$ws.Range("M136").Value = -6497.625
